$d = $word.ActiveDocument

# 1) Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2) Split the single-line mailing address into two paragraphs:
#    "2940 Sanor Pl" followed by a new paragraph "Santa Clara, CA 95051".
$addr = $d.Content
$found = $addr.Find.Execute("2940 Sanor Pl, Santa Clara CA 95051")
$addr.Text = "2940 Sanor Pl"
$addr.InsertParagraphAfter()
$cityLine = $d.Range($addr.End + 1, $addr.End + 1)
$cityLine.Text = "Santa Clara, CA 95051"

# 3) Remove the now-superfluous empty "NoSpacing" paragraph that follows
#    "Board of Directors".
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $next = $d.Paragraphs.Item($i + 1)
        $next.Range.Delete()
        break
    }
}
